$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three data rows (old rows 5,6,7) - data now fits in rows 2-4
$ws.Rows("5:7").Delete()

# Row 2 (ECs -> Cd80 -> Ctla4 -> MuSCs)
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 0.257284
$ws.Range("H2").Value = 0.771852
$ws.Range("I2").Value = 0.05106290078335718
$ws.Range("J2").Value = 0.05106290078335718
$ws.Range("M2").Value = 0.01388066666666667
$ws.Range("N2").Value = 0.041642
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.003571273442666667
$ws.Range("R2").Value = 0.032141460984
$ws.Range("S2").Value = 0.05106290078335718
$ws.Range("T2").Value = 0.05106290078335718

# Row 3 (FAPs -> Cd80 -> Ctla4 -> MuSCs)
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 4.356394000000001
$ws.Range("H3").Value = 13.069182
$ws.Range("I3").Value = 0.8646092045957484
$ws.Range("J3").Value = 0.8646092045957485
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.06046965298266668
$ws.Range("R3").Value = 0.5442268768440001
$ws.Range("S3").Value = 0.8646092045957484
$ws.Range("T3").Value = 0.8646092045957485

# Row 4 (MuSCs -> Cd80 -> Ctla4 -> MuSCs)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.424892
$ws.Range("H4").Value = 1.274676
$ws.Range("I4").Value = 0.08432789462089441
$ws.Range("J4").Value = 0.08432789462089442
$ws.Range("M4").Value = 0.01388066666666667
$ws.Range("N4").Value = 0.041642
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.005897784221333333
$ws.Range("R4").Value = 0.053080057992
$ws.Range("S4").Value = 0.08432789462089441
$ws.Range("T4").Value = 0.08432789462089442
